$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.895.13'
$ws.Range("E2").Value = '  -4.17%  '

# Row 3
$ws.Range("D3").Value = '2.447.56'
$ws.Range("E3").Value = '  -3.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.86'
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.76'
$ws.Range("E6").Value = '  -7.29%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("E7").Value = '  -4.13%  '

# Row 8
$ws.Range("E8").Value = '  +0.14%  '

# Row 9
$ws.Range("E9").Value = '  -5.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.41'
$ws.Range("E10").Value = '  -8.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0778'
$ws.Range("E11").Value = '  -3.29%  '

# Row 12
$ws.Range("E12").Value = '  -0.84%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.92'
$ws.Range("E13").Value = '  -5.80%  '

# Row 14
$ws.Range("D14").Value = '2.827.45'
$ws.Range("E14").Value = '  -3.43%  '

# Row 15
$ws.Range("D15").Value = '2.458.97'
$ws.Range("E15").Value = '  -1.88%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.37'
$ws.Range("E16").Value = '  -9.03%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  -3.71%  '

# Row 18
$ws.Range("D18").Value = '40.922.14'
$ws.Range("E18").Value = '  -4.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  -6.36%  '

# Row 20
$ws.Range("E20").Value = '  -4.38%  '

# Row 21
$ws.Range("E21").Value = '  -6.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.91'
$ws.Range("E22").Value = '  -3.67%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.57'
$ws.Range("E23").Value = '  -3.74%  '

# Row 24
$ws.Range("E24").Value = '  -5.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.91'
$ws.Range("E25").Value = '  -6.71%  '

# Row 26
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.40'
$ws.Range("E27").Value = '  -6.88%  '

# Row 28
$ws.Range("E28").Value = '  -4.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  -5.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.04'
$ws.Range("E30").Value = '  -8.66%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.62'
$ws.Range("E31").Value = '  -2.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.55'
$ws.Range("E32").Value = '  -3.95%  '

# Row 33
$ws.Range("E33").Value = '  -0.99%  '

# Row 34
$ws.Range("B34").Value = 'ApeXProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.52'
$ws.Range("E34").Value = '  -8.78%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0747'
$ws.Range("E35").Value = '  -5.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.00'
$ws.Range("E36").Value = '  -6.38%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.13'
$ws.Range("E37").Value = '  -6.65%  '

# Row 38
$ws.Range("E38").Value = '  -7.58%  '

# Row 39
$ws.Range("E39").Value = '  -4.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.103'
$ws.Range("E40").Value = '  -8.56%  '

# Row 41
$ws.Range("E41").Value = '  -4.46%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.99'
$ws.Range("E42").Value = '  -5.86%  '

# Row 43
$ws.Range("E43").Value = '  +0.10%  '

# Row 44
$ws.Range("D44").Value = '1.957.08'
$ws.Range("E44").Value = '  -1.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0283'
$ws.Range("E45").Value = '  -5.64%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  -8.77%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.69'
$ws.Range("E47").Value = '  -2.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.12'
$ws.Range("E48").Value = '  -5.89%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.04'
$ws.Range("E49").Value = '  -5.12%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '96.99'
$ws.Range("E50").Value = '  -3.91%  '

# Row 51
$ws.Range("E51").Value = '  -6.93%  '
